{"js": "// The document contains a single 20-row x 5-column table where every\n// cell holds one addition/subtraction equation (e.g. \"56+20=\"). The\n// edit replaces the text of every one of the 100 equation cells, in\n// row-major (top-to-bottom, left-to-right) order, with a new set of\n// equations while leaving all other formatting (fonts, size, paragraph\n// alignment, table layout, etc.) untouched.\n\n// New equation text, in row-major order (row 0 col 0..4, row 1 col 0..4, ...).\nconst newEquations = [\n  [\"9+65=\", \"14+26=\", \"31-24=\", \"17+14=\", \"92-12=\"],\n  [\"98-22=\", \"89-31=\", \"31+67=\", \"10+8=\", \"30-14=\"],\n  [\"35+50=\", \"64-8=\", \"69-2=\", \"34+52=\", \"65+33=\"],\n  [\"72-45=\", \"12+29=\", \"53-38=\", \"54+7=\", \"89-7=\"],\n  [\"58+41=\", \"45+29=\", \"16+28=\", \"82-25=\", \"38+46=\"],\n  [\"39-38=\", \"0+57=\", \"50+5=\", \"74-47=\", \"71-55=\"],\n  [\"57-32=\", \"79-76=\", \"84-49=\", \"23-9=\", \"32+45=\"],\n  [\"74+16=\", \"80-4=\", \"89-5=\", \"57+30=\", \"3+81=\"],\n  [\"28-21=\", \"69-25=\", \"59-44=\", \"80-5=\", \"52+22=\"],\n  [\"8+63=\", \"16+28=\", \"21-20=\", \"91-36=\", \"2+0=\"],\n  [\"72+9=\", \"84+0=\", \"8+9=\", \"41-35=\", \"73+26=\"],\n  [\"29+17=\", \"25+12=\", \"15+22=\", \"86-65=\", \"86-77=\"],\n  [\"68-6=\", \"99-29=\", \"97-70=\", \"81-1=\", \"58-2=\"],\n  [\"51-40=\", \"85-10=\", \"0+31=\", \"98-1=\", \"77+1=\"],\n  [\"71+7=\", \"56+17=\", \"15+6=\", \"42-28=\", \"41+4=\"],\n  [\"28+56=\", \"28+63=\", \"76-42=\", \"47+13=\", \"24+8=\"],\n  [\"90-64=\", \"72-11=\", \"92-12=\", \"91-70=\", \"98-26=\"],\n  [\"86-38=\", \"12+22=\", \"93-82=\", \"78-59=\", \"3+80=\"],\n  [\"59-17=\", \"24+71=\", \"54-18=\", \"15+0=\", \"25-8=\"],\n  [\"14+46=\", \"11+28=\", \"14+34=\", \"90-20=\", \"99-81=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nfor (let r = 0; r < newEquations.length && r < table.rowCount; r++) {\n  const rowValues = newEquations[r];\n  for (let c = 0; c < rowValues.length; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = rowValues[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single 20-row x 5-column table where every\n# cell holds one addition/subtraction equation (e.g. \"56+20=\"). This\n# script replaces the text of every one of the 100 equation cells, in\n# row-major (top-to-bottom, left-to-right) order, with a new set of\n# equations while leaving all other formatting (fonts, size, paragraph\n# alignment, table layout, etc.) untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# New equation text, in row-major order (row 1 col 1..5, row 2 col 1..5, ...).\n$newEquations = @(\n    @(\"9+65=\", \"14+26=\", \"31-24=\", \"17+14=\", \"92-12=\"),\n    @(\"98-22=\", \"89-31=\", \"31+67=\", \"10+8=\", \"30-14=\"),\n    @(\"35+50=\", \"64-8=\", \"69-2=\", \"34+52=\", \"65+33=\"),\n    @(\"72-45=\", \"12+29=\", \"53-38=\", \"54+7=\", \"89-7=\"),\n    @(\"58+41=\", \"45+29=\", \"16+28=\", \"82-25=\", \"38+46=\"),\n    @(\"39-38=\", \"0+57=\", \"50+5=\", \"74-47=\", \"71-55=\"),\n    @(\"57-32=\", \"79-76=\", \"84-49=\", \"23-9=\", \"32+45=\"),\n    @(\"74+16=\", \"80-4=\", \"89-5=\", \"57+30=\", \"3+81=\"),\n    @(\"28-21=\", \"69-25=\", \"59-44=\", \"80-5=\", \"52+22=\"),\n    @(\"8+63=\", \"16+28=\", \"21-20=\", \"91-36=\", \"2+0=\"),\n    @(\"72+9=\", \"84+0=\", \"8+9=\", \"41-35=\", \"73+26=\"),\n    @(\"29+17=\", \"25+12=\", \"15+22=\", \"86-65=\", \"86-77=\"),\n    @(\"68-6=\", \"99-29=\", \"97-70=\", \"81-1=\", \"58-2=\"),\n    @(\"51-40=\", \"85-10=\", \"0+31=\", \"98-1=\", \"77+1=\"),\n    @(\"71+7=\", \"56+17=\", \"15+6=\", \"42-28=\", \"41+4=\"),\n    @(\"28+56=\", \"28+63=\", \"76-42=\", \"47+13=\", \"24+8=\"),\n    @(\"90-64=\", \"72-11=\", \"92-12=\", \"91-70=\", \"98-26=\"),\n    @(\"86-38=\", \"12+22=\", \"93-82=\", \"78-59=\", \"3+80=\"),\n    @(\"59-17=\", \"24+71=\", \"54-18=\", \"15+0=\", \"25-8=\"),\n    @(\"14+46=\", \"11+28=\", \"14+34=\", \"90-20=\", \"99-81=\")\n)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    $rowValues = $newEquations[$r - 1]\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
